$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: fill in result/profit
$ws.Range("G55").Value = "Acierto"
$ws.Range("H55").Value = 0.4

# Row 58: fill in result/profit
$ws.Range("G58").Value = "Fallo"
$ws.Range("H58").Value = -1

# Row 59: fill in result/profit
$ws.Range("G59").Value = "Acierto"
$ws.Range("H59").Value = 1.62

# Rows 61-66: convert event_id (column A) from text to numeric
$ws.Range("A61").Value = 14743054
$ws.Range("A62").Value = 14743046
$ws.Range("A63").Value = 14743050
$ws.Range("A64").Value = 14743042
$ws.Range("A65").Value = 14743045
$ws.Range("A66").Value = 14743047
